# Check Date and Unit Test
# Adds a "1ER JAN" date row (B9) plus a numeric offset (C9, days) and a
# formula row (C10) that adds the two together to produce a new date.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: label, start date, number of days to add
$ws.Range("A9").Value = "1ER JAN"

# B9 and C10 share the same date format/style
$ws.Range("B9", "C10").NumberFormat = "mm-dd-yy"
$ws.Range("B9").Value = 44197

$ws.Range("C9").NumberFormat = "General"
$ws.Range("C9").Value = 42

# Row 10: compute the resulting date
$ws.Range("C10").Formula = "=B9+C9"

$ws.Range("C10").Select() | Out-Null
